$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P58").Value2 = 1981821.824
$ws.Range("T58").Value2 = 2041851.136
$ws.Range("AB58").Value2 = 1527297.792
$ws.Range("AF58").Value2 = 1875274.112

$ws.Range("H59").Value2 = -368490.048
$ws.Range("P59").Value2 = -387860.032
$ws.Range("X59").Value2 = -434668.96
$ws.Range("AF59").Value2 = -914560.896

$ws.Range("L60").Value2 = 305830.016
$ws.Range("T60").Value2 = 1418361.856
$ws.Range("X60").Value2 = 920052.032
$ws.Range("AB60").Value2 = 1060668.864
$ws.Range("AF60").Value2 = 960713.024

$ws.Range("AF61").Value2 = -326358.944

$ws.Range("AF62").Value2 = 444107.968

$ws.Range("P63").Value2 = -424255.968
$ws.Range("T63").Value2 = -196518.064
$ws.Range("AB63").Value2 = -409213.088
$ws.Range("AF63").Value2 = -166763.04

$ws.Range("L64").Value2 = -119435.024
$ws.Range("P64").Value2 = -164057.008
$ws.Range("X64").Value2 = -321834.016
$ws.Range("AB64").Value2 = -335632.928
$ws.Range("AF64").Value2 = -443946.016

$ws.Range("AF65").Value2 = -112873

$ws.Range("T66").Value2 = 85368.992
$ws.Range("AB66").Value2 = 150850.992
$ws.Range("AF66").Value2 = 167050.016

$ws.Range("P67").Value2 = -464826.048
$ws.Range("X67").Value2 = -236297.008
$ws.Range("AF67").Value2 = -213935.008

$ws.Range("H68").Value2 = -222922.016
$ws.Range("P68").Value2 = -335785.984
$ws.Range("X68").Value2 = -306926.016

$ws.Range("AB69").Value2 = 455667.936
$ws.Range("AF69").Value2 = 634353.984

$ws.Range("T73").Value2 = 499209.152
$ws.Range("AB73").Value2 = 441554.016

$ws.Range("P74").Value2 = -106651

$cols78 = @("C","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")
foreach ($col in $cols78) {
    $cell = $ws.Range($col + "78")
    $cell.Value2 = "'"
    $cell.Style = "Normal"
}

$ws.Range("T79").Value2 = 494591.136
$ws.Range("X79").Value2 = 437349.984
$ws.Range("AF79").Value2 = 346073.024

"done"
